$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row of data (row 33) following the existing pattern
$ws.Range("A33").Value = 10002
$ws.Range("B33").Value = 110032
$ws.Range("C33").Value = "eng"
$ws.Range("D33").Value = $true
$ws.Range("E33").Value = "superadmin"
$ws.Range("F33").Value = "now()"

# Update view to reflect scroll/selection position as in the edited file
$ws.Range("C31").Select()
$win = $wb.Windows.Item(1)
$win.ScrollRow = 22
$win.ScrollColumn = 1
